$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: insert $text right before $pos (a collapsed position), then
# strip bold/underline from the newly inserted run only, so it becomes
# a plain run instead of inheriting the bold+underline of the label
# that precedes it. Returns the end position of the inserted text.
# ---------------------------------------------------------------------
function Insert-PlainRun($pos, $text) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($text)
    $newRange = $d.Range($pos, $pos + $text.Length)

    # Direct Font assignment cleanly drops <w:b/> for this run.
    $newRange.Font.Bold = 0

    # Going through Find/Replace formatting-only (same text->text) cleanly
    # drops <w:u/> for this run (a plain Font.Underline=0 assignment would
    # instead write an explicit w:u w:val="none").
    $f = $newRange.Find
    $f.ClearFormatting()
    $f.Replacement.ClearFormatting()
    $f.Replacement.Font.Underline = 0
    $f.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null

    return ($pos + $text.Length)
}

# ---------------------------------------------------------------------
# 1) "Magnetische Felder:" paragraph -> append the explanatory sentence
#    as a new, non-bold / non-underlined run.
# ---------------------------------------------------------------------
$pMagnet = $d.Paragraphs.Item(36)
$labelStart = $pMagnet.Range.Start

# The paragraph mark for this paragraph used to carry bold/underline
# (matching the "Magnetische Felder:" label). Now that the paragraph will
# get plain trailing text, the mark itself should no longer be bold/
# underlined. Clearing Bold/Underline via the paragraph's own Range()
# (as opposed to a Document.Range(start,end) with identical bounds) is
# what actually reaches the paragraph-mark run-properties in pPr, not
# just the visible text runs. Do this BEFORE adding the new sentence so
# the later insertion doesn't get re-touched by this pass.
$fullRange = $pMagnet.Range()
$fullRange.Font.Bold = 0
$fullRange.Font.Underline = 0

$labelRange = $d.Range($labelStart, $labelStart + 19)  # "Magnetische Felder:" = 19 chars
Write-Output ("label text=[" + $labelRange.Text + "]")
$labelRange.Font.Bold = -1
$labelRange.Font.Underline = 1

$insertPos = $pMagnet.Range.End - 1
$sentence = " Wie ein elektrischen Feld gibt es auch ein magnetisches Feld in der Physik. Diese ist gekennzeichnet durch einen Nord- und einen Südpol, welche nicht trennbar sind. Ein solches magnetisches Feld kann durch Stromfluss erzeugt werden."
Insert-PlainRun $insertPos $sentence | Out-Null

Write-Output "step1 done"
